$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-26 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-27 Thursday", 2)
$d.Content.Find.Execute("159÷5=31, 4", $true, $false, $false, $false, $false, $true, 1, $false, "230÷9=25, 5", 2)
$d.Content.Find.Execute("464÷8=58, 0", $true, $false, $false, $false, $false, $true, 1, $false, "316÷2=158, 0", 2)
$d.Content.Find.Execute("687÷3=229, 0", $true, $false, $false, $false, $false, $true, 1, $false, "266÷9=29, 5", 2)
$d.Content.Find.Execute("573÷7=81, 6", $true, $false, $false, $false, $false, $true, 1, $false, "856÷7=122, 2", 2)
$d.Content.Find.Execute("334÷3=111, 1", $true, $false, $false, $false, $false, $true, 1, $false, "233÷8=29, 1", 2)
$d.Content.Find.Execute("646÷3=215, 1", $true, $false, $false, $false, $false, $true, 1, $false, "910÷3=303, 1", 2)
$d.Content.Find.Execute("878÷8=109, 6", $true, $false, $false, $false, $false, $true, 1, $false, "558÷9=62, 0", 2)
$d.Content.Find.Execute("823÷8=102, 7", $true, $false, $false, $false, $false, $true, 1, $false, "421÷8=52, 5", 2)
$d.Content.Find.Execute("666÷4=166, 2", $true, $false, $false, $false, $false, $true, 1, $false, "194÷3=64, 2", 2)
$d.Content.Find.Execute("366÷8=45, 6", $true, $false, $false, $false, $false, $true, 1, $false, "782÷6=130, 2", 2)
$d.Content.Find.Execute("828÷6=138, 0", $true, $false, $false, $false, $false, $true, 1, $false, "222÷2=111, 0", 2)
$d.Content.Find.Execute("209÷4=52, 1", $true, $false, $false, $false, $false, $true, 1, $false, "582÷6=97, 0", 2)
$d.Content.Find.Execute("971÷2=485, 1", $true, $false, $false, $false, $false, $true, 1, $false, "522÷3=174, 0", 2)
$d.Content.Find.Execute("161÷9=17, 8", $true, $false, $false, $false, $false, $true, 1, $false, "108÷9=12, 0", 2)
$d.Content.Find.Execute("982÷7=140, 2", $true, $false, $false, $false, $false, $true, 1, $false, "280÷9=31, 1", 2)
$d.Content.Find.Execute("589÷7=84, 1", $true, $false, $false, $false, $false, $true, 1, $false, "759÷8=94, 7", 2)
$d.Content.Find.Execute("158÷6=26, 2", $true, $false, $false, $false, $false, $true, 1, $false, "334÷5=66, 4", 2)
$d.Content.Find.Execute("705÷8=88, 1", $true, $false, $false, $false, $false, $true, 1, $false, "871÷6=145, 1", 2)
$d.Content.Find.Execute("490÷7=70, 0", $true, $false, $false, $false, $false, $true, 1, $false, "133÷7=19, 0", 2)
$d.Content.Find.Execute("707÷7=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "964÷4=241, 0", 2)
$d.Content.Find.Execute("141÷5=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "123÷9=13, 6", 2)
$d.Content.Find.Execute("923÷5=184, 3", $true, $false, $false, $false, $false, $true, 1, $false, "483÷3=161, 0", 2)
$d.Content.Find.Execute("686÷3=228, 2", $true, $false, $false, $false, $false, $true, 1, $false, "391÷8=48, 7", 2)
$d.Content.Find.Execute("676÷5=135, 1", $true, $false, $false, $false, $false, $true, 1, $false, "346÷4=86, 2", 2)
$d.Content.Find.Execute("627÷9=69, 6", $true, $false, $false, $false, $false, $true, 1, $false, "629÷7=89, 6", 2)
